$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted right before the existing
# row 264 ("Femacal de La Calera" / Pepino ensalada data block). Insert a
# fresh row there so every following row shifts down by one (264->265,
# ..., 280->281), matching the diff exactly.
$ws.Rows.Item(264).Insert()

# Populate the newly inserted row with the new observation's data. The
# static/categorical columns repeat the same values used by every other
# row in this data block.
$ws.Cells.Item(264, 1).Value = 3
$ws.Cells.Item(264, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(264, 3).Value = "Coquimbo"
$ws.Cells.Item(264, 4).Value = 44585
$ws.Cells.Item(264, 5).Value = 5
$ws.Cells.Item(264, 6).Value = 100112043
$ws.Cells.Item(264, 7).Value = "Pepino ensalada"
$ws.Cells.Item(264, 8).Value = "Sin especificar"
$ws.Cells.Item(264, 9).Value = "Primera"
$ws.Cells.Item(264, 10).Value = 100
$ws.Cells.Item(264, 11).Value = 11000
$ws.Cells.Item(264, 12).Value = 12000
$ws.Cells.Item(264, 13).Value = 11500
$ws.Cells.Item(264, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(264, 15).Value = "Limache"
$ws.Cells.Item(264, 16).Value = 164
$ws.Cells.Item(264, 17).Value = 70
$ws.Cells.Item(264, 18).Value = "Hortaliza"
